# Updates the cryptos price/volume table (columns D and E) on the active
# worksheet to reflect refreshed market data, mirroring the values pulled
# in by the "Updated cryptos list" GitHub Actions workflow run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes a single cell update: the row number, the target
# column letter (D = Price, E = Volume(1h)) and the new text value.
$updates = @(
    @{Row=2; Col="D"; Value="90.767.74"},
    @{Row=2; Col="E"; Value="  +2.48%  "},
    @{Row=3; Col="D"; Value="3.173.62"},
    @{Row=3; Col="E"; Value="  -0.56%  "},
    @{Row=4; Col="E"; Value="  -0.04%  "},
    @{Row=5; Col="D"; Value="215.32"},
    @{Row=5; Col="E"; Value="  +2.70%  "},
    @{Row=6; Col="D"; Value="631.36"},
    @{Row=6; Col="E"; Value="  +2.71%  "},
    @{Row=7; Col="D"; Value="0.398"},
    @{Row=7; Col="E"; Value="  +1.81%  "},
    @{Row=8; Col="D"; Value="0.726"},
    @{Row=8; Col="E"; Value="  +5.65%  "},
    @{Row=9; Col="D"; Value="1.00"},
    @{Row=10; Col="D"; Value="3.168.84"},
    @{Row=10; Col="E"; Value="  -0.50%  "},
    @{Row=11; Col="D"; Value="0.568"},
    @{Row=11; Col="E"; Value="  +3.58%  "},
    @{Row=12; Col="D"; Value="0.182"},
    @{Row=12; Col="E"; Value="  +2.83%  "},
    @{Row=13; Col="E"; Value="  +1.27%  "},
    @{Row=14; Col="D"; Value="90.673.79"},
    @{Row=14; Col="E"; Value="  +2.76%  "},
    @{Row=15; Col="D"; Value="5.33"},
    @{Row=15; Col="E"; Value="  -0.20%  "},
    @{Row=16; Col="D"; Value="3.767.75"},
    @{Row=16; Col="E"; Value="  -0.20%  "},
    @{Row=17; Col="D"; Value="32.54"},
    @{Row=17; Col="E"; Value="  -0.78%  "},
    @{Row=18; Col="D"; Value="3.190.74"},
    @{Row=18; Col="E"; Value="  +0.78%  "},
    @{Row=19; Col="E"; Value="  +4.03%  "},
    @{Row=20; Col="D"; Value="0.0000213"},
    @{Row=20; Col="E"; Value="  +35.18%  "},
    @{Row=21; Col="D"; Value="13.41"},
    @{Row=21; Col="E"; Value="  -1.62%  "},
    @{Row=22; Col="D"; Value="433.97"},
    @{Row=22; Col="E"; Value="  +3.84%  "},
    @{Row=23; Col="D"; Value="8.45"},
    @{Row=23; Col="E"; Value="  -1.19%  "},
    @{Row=24; Col="D"; Value="4.98"},
    @{Row=24; Col="E"; Value="  -3.57%  "},
    @{Row=25; Col="D"; Value="5.30"},
    @{Row=25; Col="E"; Value="  -1.28%  "},
    @{Row=26; Col="D"; Value="11.65"},
    @{Row=26; Col="E"; Value="  -6.55%  "},
    @{Row=27; Col="D"; Value="80.81"},
    @{Row=27; Col="E"; Value="  +9.14%  "},
    @{Row=28; Col="D"; Value="3.346.49"},
    @{Row=28; Col="E"; Value="  +0.39%  "},
    @{Row=29; Col="E"; Value="  +0.11%  "},
    @{Row=30; Col="D"; Value="0.160"},
    @{Row=30; Col="E"; Value="  -5.27%  "},
    @{Row=31; Col="E"; Value="  +0.14%  "},
    @{Row=32; Col="D"; Value="4.01"},
    @{Row=32; Col="E"; Value="  +26.21%  "},
    @{Row=33; Col="D"; Value="8.34"},
    @{Row=33; Col="E"; Value="  -0.32%  "},
    @{Row=34; Col="D"; Value="512.08"},
    @{Row=34; Col="E"; Value="  -7.61%  "},
    @{Row=35; Col="D"; Value="6.96"},
    @{Row=35; Col="E"; Value="  -0.53%  "},
    @{Row=36; Col="D"; Value="1.88"},
    @{Row=36; Col="E"; Value="  +0.42%  "},
    @{Row=37; Col="D"; Value="1.28"},
    @{Row=37; Col="E"; Value="  -3.47%  "},
    @{Row=38; Col="D"; Value="22.27"},
    @{Row=38; Col="E"; Value="  +0.70%  "},
    @{Row=39; Col="D"; Value="22.37"},
    @{Row=39; Col="E"; Value="  +2.53%  "},
    @{Row=40; Col="E"; Value="  +0.28%  "},
    @{Row=41; Col="D"; Value="0.126"},
    @{Row=41; Col="E"; Value="  -4.05%  "},
    @{Row=42; Col="E"; Value="  -0.04%  "},
    @{Row=43; Col="D"; Value="1.91"},
    @{Row=43; Col="E"; Value="  -1.54%  "},
    @{Row=44; Col="D"; Value="0.372"},
    @{Row=44; Col="E"; Value="  -1.80%  "},
    @{Row=45; Col="D"; Value="146.22"},
    @{Row=45; Col="E"; Value="  -2.80%  "},
    @{Row=46; Col="D"; Value="43.84"},
    @{Row=46; Col="E"; Value="  +1.24%  "},
    @{Row=47; Col="D"; Value="168.99"},
    @{Row=47; Col="E"; Value="  -4.02%  "},
    @{Row=48; Col="D"; Value="0.125"},
    @{Row=48; Col="E"; Value="  -1.01%  "},
    @{Row=49; Col="D"; Value="0.735"},
    @{Row=49; Col="E"; Value="  +5.18%  "},
    @{Row=50; Col="D"; Value="24.76"},
    @{Row=50; Col="E"; Value="  +0.59%  "},
    @{Row=51; Col="E"; Value="  -3.98%  "}
)

foreach ($u in $updates) {
    $cellRef = "$($u.Col)$($u.Row)"
    $range = $ws.Range($cellRef)
    $value = $u.Value

    if ($u.Col -eq "D") {
        # Price values are plain text in this sheet (e.g. "90.767.74" or
        # "1.00"). Some of the new prices parse as plain numbers (e.g.
        # "215.32"), and Excel would silently convert those to numeric
        # values - dropping significant trailing zeros - unless we force
        # them to stay text by using a leading quote prefix, exactly as
        # typing an apostrophe before a number does in the Excel UI.
        # Values containing more than one "." (like "90.767.74") already
        # cannot be parsed as numbers, so they do not need the prefix.
        $dotCount = ($value.ToCharArray() | Where-Object { $_ -eq '.' }).Count
        $looksNumeric = $dotCount -le 1
        if ($looksNumeric) {
            $range.Value = "'" + $value
        } else {
            $range.Value = $value
        }
    } else {
        # Volume percentages always contain a trailing "%" plus padding
        # spaces, so Excel already keeps these as plain text.
        $range.Value = $value
    }
}
